# planety_quizz.xlsx - implementacia Toda a zacatie doplnania planet
# Replaces the Slnko "weight / gravity" questions (rows 7-8) with new
# Slnko questions (corona/sunspots) and replaces the placeholder
# "Otazka Merkur 1..10" rows (rows 9-18) with real Merkur questions.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 7: new Slnko question (was "Ake tazke je Slnko?") ---
$ws.Range("B7").Value = "Aká je najvrchnejšia časť Slnka?"
$ws.Range("E7").Value = "Koróna"
$ws.Range("C7").Value = "Fotosféra"
$ws.Range("D7").Value = "Niečo žlté"

# --- Row 8: new Slnko question (was "Co robi gravitacne pole Slnka?") ---
$ws.Range("B8").Value = "Čo sú slnečné škvrny?"
$ws.Range("E8").Value = "Chladnejšie časti Slnka"
$ws.Range("C8").Value = "Slnko sa neumylo"
$ws.Range("D8").Value = "Niekto tam nechal plachtu"

# --- Row 9: "Kde sa nachadza Merkur?" (was "Otazka Merkur 1") ---
$ws.Range("B9").Value = "Kde sa nachádza Merkúr?"
$ws.Range("E9").Value = "Najbližšie ku Slnku"
$ws.Range("C9").Value = "Najďalej pri Slnku"
$ws.Range("D9").Value = "V n124qw galaxií"

# --- Row 10: "Aka je velkost Merkuru" (was "Otazka Merkur 2") ---
$ws.Range("B10").Value = "Aká je veľkosť Merkúru"
$ws.Range("E10").Value = "Patrí medzi druhú najmenšiu"
$ws.Range("C10").Value = "Je najväčšia"
$ws.Range("D10").Value = "Skoro rovnaká ako Zem"

# --- Row 11: "Ma nejaku atmosferu?" (was "Otazka Merkur 3") ---
$ws.Range("B11").Value = "Má nejakú atmosféru?"
$ws.Range("E11").Value = "Takmer žiadnu"
$ws.Range("D11").Value = "Čo je armosféra?"
$ws.Range("C11").Value = "Áno"

# --- Row 12: "Ake su najvacsie teploty?" (was "Otazka Merkur 4") ---
$ws.Range("B12").Value = "Aké sú najväčšie teploty?"
$ws.Range("E12").Value = "430 C"
$ws.Range("D12").Value = "1500 C"
$ws.Range("C12").Value = "589 C"

# --- Row 13: "Aka je najnizsia teplota?" (was "Otazka Merkur 5") ---
$ws.Range("B13").Value = "Aká je najnižšia teplota?"
$ws.Range("C13").Value = 55
$ws.Range("D13").Value = 30
$ws.Range("E13").Value = -180

# --- Row 14: "Cim je pokryty Merkur?" (was "Otazka Merkur 6") ---
$ws.Range("B14").Value = "Čím je pokrytý Merkúr?"

# --- Row 15: "Kolko mesiacov ma Merkur?" (was "Otazka Merkur 7") ---
$ws.Range("B15").Value = "Koľko mesiacov má Merkúr?"
$ws.Range("C15").Value = 1
$ws.Range("D15").Value = 4
$ws.Range("E15").Value = 0

# --- Row 16: "Kolko trva dlzka roka na Merkure?" (was "Otazka Merkur 8") ---
$ws.Range("B16").Value = "Koľko trvá dĺžka roka  na Merkúre?"
$ws.Range("C16").Value = 66
$ws.Range("D16").Value = 77
$ws.Range("E16").Value = 88

# --- Row 17: "Preco nema takmer ziadnu atmosferu?" (was "Otazka Merkur 9") ---
$ws.Range("B17").Value = "Prečo nemá takmer žiadnu atmosféru?"
$ws.Range("E17").Value = "Hmotnosť planéty je príliš malá"
$ws.Range("D17").Value = "Atmosféra nemá rada Merkúr"
$ws.Range("C17").Value = "Čo je to atmosféra?"

# --- Row 18: "Kolko rocnych obdobi ma Merkur?" (was "Otazka Merkur 10") ---
$ws.Range("B18").Value = "Koľko ročných období má Merkúr?"
$ws.Range("E18").Value = "Žiadne"
$ws.Range("D18").Value = "Jar a Leto"
$ws.Range("C18").Value = 4

# Update selection to match the authored edit (cursor ended on B18)
$ws.Range("B18").Select()
